# Generate Report for Handoff
# Updates the "b.md" row (row 3) on all three sheets to reflect that the
# handoff package for file "b" is now ready / has been handed off, mirroring
# the structure already present for file "a".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": zh-cn / de-de status + handoff date columns for b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-31 07:02:13"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Status / Latest Handoff File / Latest Handoff Datetime
# for the b.md row
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-31 07:02:03"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de": Status / Latest Handoff File / Latest Handoff Datetime
# for the b.md row
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-31 07:02:13"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
